$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.486320018768311
$ws.Range("B1").Value = 4.374085426330566
$ws.Range("C1").Value = 2.938593626022339
$ws.Range("D1").Value = 1.934455156326294
$ws.Range("E1").Value = 1.291341543197632
